$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $r = $ws.Range($addr)
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        # Value looks like a plain number (e.g. "10.50"); force text
        # storage so Excel does not silently coerce it to a Double
        # and drop significant trailing zeros / formatting.
        $r.NumberFormat = "@"
        $r.Value = $val
        $r.ClearFormats()
    } else {
        $r.Value = $val
    }
}

$changes = @(
    @{ Cell = "D2"; Value = '70.974.13' },
    @{ Cell = "E2"; Value = '  +0.55%  ' },
    @{ Cell = "D3"; Value = '3.844.33' },
    @{ Cell = "E3"; Value = '  +1.10%  ' },
    @{ Cell = "E4"; Value = '  -0.03%  ' },
    @{ Cell = "D5"; Value = '689.45' },
    @{ Cell = "E5"; Value = '  +3.33%  ' },
    @{ Cell = "D6"; Value = '172.69' },
    @{ Cell = "E6"; Value = '  +2.12%  ' },
    @{ Cell = "D7"; Value = '3.841.57' },
    @{ Cell = "E7"; Value = '  +1.07%  ' },
    @{ Cell = "E8"; Value = '  +0.01%  ' },
    @{ Cell = "E10"; Value = '  +1.40%  ' },
    @{ Cell = "D11"; Value = '7.41' },
    @{ Cell = "E11"; Value = '  +4.77%  ' },
    @{ Cell = "E12"; Value = '  -0.48%  ' },
    @{ Cell = "E13"; Value = '  +5.70%  ' },
    @{ Cell = "D14"; Value = '36.57' },
    @{ Cell = "E14"; Value = '  +2.31%  ' },
    @{ Cell = "D15"; Value = '4.493.03' },
    @{ Cell = "E15"; Value = '  +1.13%  ' },
    @{ Cell = "D16"; Value = '3.857.54' },
    @{ Cell = "E16"; Value = '  +1.46%  ' },
    @{ Cell = "D17"; Value = '71.019.08' },
    @{ Cell = "E17"; Value = '  +0.74%  ' },
    @{ Cell = "D18"; Value = '17.77' },
    @{ Cell = "E18"; Value = '  +0.50%  ' },
    @{ Cell = "E19"; Value = '  +0.50%  ' },
    @{ Cell = "E20"; Value = '  +0.31%  ' },
    @{ Cell = "D21"; Value = '11.09' },
    @{ Cell = "E21"; Value = '  -4.36%  ' },
    @{ Cell = "D22"; Value = '487.19' },
    @{ Cell = "E22"; Value = '  +2.73%  ' },
    @{ Cell = "E23"; Value = '  +0.89%  ' },
    @{ Cell = "D24"; Value = '84.69' },
    @{ Cell = "E24"; Value = '  +2.08%  ' },
    @{ Cell = "E25"; Value = '  +2.31%  ' },
    @{ Cell = "D26"; Value = '12.35' },
    @{ Cell = "E26"; Value = '  +1.17%  ' },
    @{ Cell = "D27"; Value = '10.50' },
    @{ Cell = "E27"; Value = '  +1.46%  ' },
    @{ Cell = "E28"; Value = '  +1.00%  ' },
    @{ Cell = "D29"; Value = '3.998.26' },
    @{ Cell = "E29"; Value = '  +1.14%  ' },
    @{ Cell = "E30"; Value = '  +0.11%  ' },
    @{ Cell = "E31"; Value = '  +9.14%  ' },
    @{ Cell = "D32"; Value = '7.62' },
    @{ Cell = "E32"; Value = '  +2.75%  ' },
    @{ Cell = "E33"; Value = '  +0.16%  ' },
    @{ Cell = "D34"; Value = '29.71' },
    @{ Cell = "E34"; Value = '  +0.57%  ' },
    @{ Cell = "E35"; Value = '  +2.46%  ' },
    @{ Cell = "D37"; Value = '3.797.07' },
    @{ Cell = "E37"; Value = '  +1.02%  ' },
    @{ Cell = "E38"; Value = '  -0.06%  ' },
    @{ Cell = "E39"; Value = '  +1.05%  ' },
    @{ Cell = "E40"; Value = '  +12.83%  ' },
    @{ Cell = "E41"; Value = '  +0.40%  ' },
    @{ Cell = "D42"; Value = '6.05' },
    @{ Cell = "E42"; Value = '  +1.63%  ' },
    @{ Cell = "E43"; Value = '  +5.00%  ' },
    @{ Cell = "D44"; Value = '0.999' },
    @{ Cell = "E44"; Value = '  -0.02%  ' },
    @{ Cell = "E45"; Value = '  +0.05%  ' },
    @{ Cell = "D46"; Value = '164.69' },
    @{ Cell = "E46"; Value = '  +3.77%  ' },
    @{ Cell = "D47"; Value = '0.000307' },
    @{ Cell = "E47"; Value = '  +7.21%  ' },
    @{ Cell = "E48"; Value = '  +1.35%  ' },
    @{ Cell = "D49"; Value = '44.43' },
    @{ Cell = "E49"; Value = '  -2.63%  ' },
    @{ Cell = "E50"; Value = '  +1.11%  ' },
    @{ Cell = "E51"; Value = '  -2.86%  ' }
)

foreach ($c in $changes) {
    Set-CellText $c.Cell $c.Value
}

Write-Host "Updated $($changes.Count) cells"
